# Regenerate the whole experiment data on Arkusz1 (sheet1):
# new graph_id count (0..32), updated inside_prob / outside_prob numbers
# (now stored as plain numbers, not shared-string lookups), 4 new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$data = @(
    @(0, 125, 2, 0.4, 0.01),
    @(1, 125, 2, 0.3, 0.02),
    @(2, 125, 2, 0.5, 0.05),
    @(3, 125, 2, 0.4, 0.08),
    @(4, 125, 2, 0.4, 0.1),
    @(5, 125, 2, 0.7, 0.12),
    @(6, 125, 3, 0.35, 0.06),
    @(7, 125, 3, 0.3, 0.02),
    @(8, 125, 3, 0.4, 0.05),
    @(9, 125, 3, 0.4, 0.02),
    @(10, 125, 3, 0.4, 0.05),
    @(11, 125, 3, 0.4, 0.08),
    @(12, 125, 4, 0.4, 0.01),
    @(13, 125, 4, 0.3, 0.03),
    @(14, 125, 4, 0.5, 0.02),
    @(15, 125, 4, 0.3, 0.02),
    @(16, 125, 4, 0.3, 0.02),
    @(17, 125, 4, 0.35, 0.02),
    @(18, 125, 4, 0.5, 0.07),
    @(19, 125, 5, 0.5, 0.05),
    @(20, 125, 5, 0.4, 0.01),
    @(21, 125, 5, 0.3, 0.02),
    @(22, 125, 5, 0.3, 0.02),
    @(23, 125, 5, 0.4, 0.05),
    @(24, 125, 5, 0.5, 0.08),
    @(25, 125, 6, 0.45, 0.1),
    @(26, 125, 6, 0.4, 0.01),
    @(27, 125, 6, 0.5, 0.02),
    @(28, 125, 6, 0.4, 0.05),
    @(29, 125, 6, 0.4, 0.03),
    @(30, 125, 6, 0.8, 0.06),
    @(31, 125, 6, 0.65, 0.08),
    @(32, 125, 6, 0.6, 0.05)
)

$firstDataRow = 2
$lastOldRow = 30
$lastNewRow = $firstDataRow + $data.Length - 1

# Clear out the old data block first (old D/E columns held shared-string
# text values; new ones are plain numbers so a clean overwrite avoids any
# stale type hanging around), then clear any now-unused trailing rows.
$ws.Range("A$firstDataRow`:E$lastOldRow").Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $firstDataRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

# Selection / view now anchors on column A.
[void]$ws.Range("A2:A$lastNewRow").Select()
$excel.ActiveWindow.Zoom = 100
